$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell G1 ("sum") onto the new
# header cell H1, then set its value to "Save" (matches the bold/centered/
# bordered header style, i.e. cellXfs style index 1, used by the other
# headers).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for the "Save" column, numeric value 0.
$ws.Range("H2").Value = 0
